$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/3/2025  Through  2/9/2025"

# --- Row 16 (Robbery): simple numeric overwrites, styles unchanged ---
$ws.Range("D16").Value = 2
$ws.Range("G16").Value = 3
$ws.Range("J16").Value = 3

# --- Row 17 (Fel. Assault): C17 numeric(2) -> text "0" (same state as D17) ---
$ws.Range("D17").Copy($ws.Range("C17"))

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 2

# D19: text "0" -> numeric 1, style13 -> style14
$ws.Range("C19").Copy($ws.Range("D19"))
$ws.Range("D19").Value = 1

# E19: text " " -> numeric 100, style13 -> style15
$ws.Range("H19").Copy($ws.Range("E19"))
$ws.Range("E19").Value = 100

$ws.Range("F19").Value = 3
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 50
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 50
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -25
$ws.Range("N19").Value = -57.142857142857

# --- Row 21 (TOTAL): simple numeric overwrites, styles unchanged ---
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 5
$ws.Range("J21").Value = 6
$ws.Range("K21").Value = -16.666666666666
$ws.Range("L21").Value = 25
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = -82.142857142857

# --- Row 24 (Petit Larceny) ---
$ws.Range("L24").Value = -80

# --- Row 26 (Misd. Assault): L26 text " " -> numeric -50, style13 -> style15 ---
$ws.Range("K26").Copy($ws.Range("L26"))
$ws.Range("L26").Value = -50

# --- Row 28 (Other Sex Crimes): C28/F28/I28 text "0" -> numeric 1, style13 -> style14 ---
$ws.Range("C19").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1

$ws.Range("C19").Copy($ws.Range("F28"))
$ws.Range("F28").Value = 1

$ws.Range("C19").Copy($ws.Range("I28"))
$ws.Range("I28").Value = 1

$ws.Range("L28").Value = 0
